$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.290.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "'1.815.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.31%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'329.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4421"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.82%  "
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").Value = "'44.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'0.07704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'22.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'6.256"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").Value = "'7.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("D16").Value = "'1.819.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("D17").Value = "'92.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.94%  "
$ws.Range("D18").Value = "'0.00001083"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").Value = "'0.06645"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.42%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("E21").Value = "  +3.87%  "
$ws.Range("D22").Value = "'6.209"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").Value = "'28.345.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("D24").Value = "'11.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.56%  "
$ws.Range("D25").Value = "'1.994"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -16.99%  "
$ws.Range("D26").Value = "'20.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("D27").Value = "'155.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.16%  "
$ws.Range("D28").Value = "'2.023.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "'128.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Value = "'1.206"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.99%  "
$ws.Range("D32").Value = "'5.866"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").Value = "'0.09212"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "'3.659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'13.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("D36").Value = "'0.02356"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.31%  "
$ws.Range("D37").Value = "'0.2170"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06223"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.163"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6576"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").Value = "'1.197"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "'8.144"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").Value = "'0.9990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "'13.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'1.386"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("D46").Value = "'0.6082"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("D47").Value = "'3.762"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.038"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.94%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'126.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "'1.154"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.51%  "
$ws.Range("D51").Value = "'0.06981"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.26%  "
